# Apply cryptocurrency price/volume refresh to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price column (D) updates ---
# Values such as "1.001" would otherwise be auto-parsed by Excel as numbers,
# so we temporarily force Text format, assign the literal string, then clear
# the temporary formatting override to leave the cell style untouched.
$priceUpdates = [ordered]@{
    'D2' = '28.185.88'
    'D3' = '1.779.64'
    'D4' = '1.001'
    'D5' = '243.74'
    'D6' = '1.001'
    'D7' = '0.4906'
    'D8' = '0.2662'
    'D9' = '0.06245'
    'D10' = '1.777.71'
    'D11' = '16.46'
    'D12' = '0.07016'
    'D13' = '0.6260'
    'D14' = '4.629'
    'D15' = '79.82'
    'D16' = '28.181.02'
    'D17' = '1.0000'
    'D18' = '0.9998'
    'D19' = '0.000007229'
    'D21' = '2.004.87'
    'D23' = '8.723'
    'D24' = '5.220'
    'D25' = '140.87'
    'D26' = '15.75'
    'D27' = '1.853'
    'D28' = '108.96'
    'D29' = '1.401'
    'D30' = '4.168'
    'D31' = '0.08239'
    'D32' = '3.757'
    'D33' = '0.04887'
    'D34' = '1.068'
    'D35' = '2.611'
    'D36' = '0.6486'
    'D37' = '0.9476'
    'D38' = '2.595'
    'D39' = '2.045'
    'D40' = '5.896'
    'D41' = '0.01545'
    'D42' = '0.9994'
    'D43' = '99.32'
    'D44' = '0.3968'
    'D45' = '7.146'
    'D47' = '0.05434'
    'D48' = '8.011'
    'D49' = '1.294'
    'D50' = '30.60'
    'D51' = '52.76'
}
foreach ($addr in $priceUpdates.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $priceUpdates[$addr]
    $cell.ClearFormats()
}

# --- Volume(1h) column (E) updates ---
$volumeUpdates = [ordered]@{
    'E2' = '  +5.67%  '
    'E3' = '  +3.07%  '
    'E4' = '  +0.33%  '
    'E5' = '  +0.83%  '
    'E6' = '  +0.25%  '
    'E7' = '  -0.48%  '
    'E8' = '  +2.07%  '
    'E9' = '  +0.56%  '
    'E10' = '  +3.02%  '
    'E11' = '  +4.06%  '
    'E12' = '  +0.57%  '
    'E13' = '  +2.49%  '
    'E14' = '  +3.08%  '
    'E15' = '  +3.42%  '
    'E16' = '  +6.37%  '
    'E17' = '  +0.16%  '
    'E18' = '  +0.19%  '
    'E19' = '  +0.18%  '
    'E20' = '  +5.97%  '
    'E21' = '  +2.71%  '
    'E22' = '  +2.19%  '
    'E23' = '  +2.06%  '
    'E24' = '  +2.72%  '
    'E25' = '  +2.29%  '
    'E26' = '  +2.80%  '
    'E27' = '  +4.81%  '
    'E28' = '  +2.44%  '
    'E29' = '  +1.03%  '
    'E30' = '  +6.40%  '
    'E31' = '  +3.39%  '
    'E32' = '  +2.25%  '
    'E33' = '  +9.18%  '
    'E34' = '  +6.70%  '
    'E35' = '  -0.05%  '
    'E36' = '  +4.00%  '
    'E37' = '  +1.22%  '
    'E38' = '  +7.55%  '
    'E39' = '  +0.30%  '
    'E40' = '  +5.40%  '
    'E41' = '  +2.45%  '
    'E42' = '  -0.11%  '
    'E43' = '  -0.06%  '
    'E44' = '  +3.18%  '
    'E45' = '  +3.97%  '
    'E46' = '  +4.09%  '
    'E47' = '  +1.04%  '
    'E48' = '  +2.27%  '
    'E49' = '  +5.35%  '
    'E50' = '  +1.36%  '
    'E51' = '  +2.31%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
